$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: only the observation id (B19) changes ---
$ws.Range("B19").Value = 91805

# --- Rows 20 and 21: swap the two observation records ---
# (row 20 becomes what row 21 was, and vice versa - except B20 is a genuinely
#  new id, 79244, rather than the old B21 value of 79243)

$ws.Range("A20").Value = 131090145
$ws.Range("B20").Value = 79244
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("M20").Value = ""
$ws.Range("P20").Value = "Timmeråsen, Timmeråsen, Ång"
$ws.Range("Q20").Value = 585013
$ws.Range("R20").Value = 7060142
$ws.Range("S20").Value = 10
$ws.Range("Z20").Value = "13:53"
$ws.Range("AB20").Value = "13:53"
$ws.Range("AC20").Value = ""
$ws.Range("AW20").Value = "Kim Hultgren"
$ws.Range("AX20").Value = "Kim Hultgren"

$ws.Range("A21").Value = 131090091
$ws.Range("B21").Value = 57884
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("M21").Value = "färska spår"
$ws.Range("P21").Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Range("Q21").Value = 585024
$ws.Range("R21").Value = 7060099
$ws.Range("S21").Value = 15
$ws.Range("Z21").Value = "13:50"
$ws.Range("AB21").Value = "13:50"
$ws.Range("AC21").Value = "Färska ringhack, tall"
$ws.Range("AW21").Value = "Daniel Rutschman"
$ws.Range("AX21").Value = "Daniel Rutschman"

# --- Rows 32 and 33: swap the two observation records exactly ---

$ws.Range("A32").Value = 131144496
$ws.Range("B32").Value = 57884
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = "Tretåig hackspett"
$ws.Range("G32").Value = "Picoides tridactylus"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("M32").Value = "färska spår"
$ws.Range("Q32").Value = 584875
$ws.Range("R32").Value = 7060422
$ws.Range("AC32").Value = "Färska ringhack, tall"

$ws.Range("A33").Value = 131144498
$ws.Range("B33").Value = 58043
$ws.Range("E33").Value = 103021
$ws.Range("F33").Value = "Talltita"
$ws.Range("G33").Value = "Poecile montanus"
$ws.Range("H33").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M33").Value = "lockläte, övriga läten"
$ws.Range("Q33").Value = 584857
$ws.Range("R33").Value = 7060494
$ws.Range("AC33").Value = ""
